$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# E5: expired days count increments by one
$ws.Range("E5").Value = "❌ EXPIRED 3761 days ago"

# E6: days-to-expiry decrements by one
$ws.Range("E6").Value = "⚠️ Expires in 18 days"

# B7: certificate expiration date pushed out to 2026-07-17 (keep as text, like the rest of column B)
$b7 = $ws.Range("B7")
$b7.NumberFormat = "@"
$b7.Value = "2026-07-17"
$b7.ClearFormats()

# E7: no longer expiring soon (cert was renewed), status cleared
$ws.Range("E7").Value = ""

# E8: days-to-expiry decrements by one
$ws.Range("E8").Value = "⚠️ Expires in 18 days"
